$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 883.3333
$ws.Range("I80").Value = 1100
$ws.Range("J80").Value = 775
$ws.Range("K80").Value = 3300
$ws.Range("L80").Value = 2325
$ws.Range("M80").Value = -2302
$ws.Range("N80").Value = -4321
$ws.Range("H83").Value = 883.3333
$ws.Range("I83").Value = 1100
$ws.Range("J83").Value = 775
$ws.Range("K83").Value = 9900
$ws.Range("L83").Value = 6975
$ws.Range("M83").Value = -4908
$ws.Range("N83").Value = -16959

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2817.5557
$ws.Range("I97").Value = 2051.2856
$ws.Range("J97").Value = 5499.5
$ws.Range("K97").Value = 2051.2856
$ws.Range("L97").Value = 5499.5
$ws.Range("M97").Value = -1555.2856
$ws.Range("N97").Value = -6491.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 721
$ws.Range("I94").Value = 774.3333
$ws.Range("J94").Value = 614.3333
$ws.Range("K94").Value = 774.3333
$ws.Range("L94").Value = 614.3333
$ws.Range("M94").Value = -323.3333
$ws.Range("N94").Value = -1516.3333
$ws.Range("H99").Value = 1499.5
$ws.Range("I99").Value = 1499.5
$ws.Range("K99").Value = 1499.5
$ws.Range("M99").Value = -1.5
$ws.Range("H140").Value = 95779.25
$ws.Range("J140").Value = 95779.25
$ws.Range("L140").Value = 95779.25
$ws.Range("N140").Value = -106139.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2680.5454
$ws.Range("I134").Value = 2667.2
$ws.Range("K134").Value = 8001.599999999999
$ws.Range("M134").Value = -5466.599999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1000
$ws.Range("J114").Value = 1000
$ws.Range("L114").Value = 3000
$ws.Range("N114").Value = -9508
$ws.Range("H140").Value = 4000
$ws.Range("J140").Value = 4000
$ws.Range("L140").Value = 12000
$ws.Range("N140").Value = -22360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 103
$ws.Range("N2").ClearContents()
$ws.Range("H80").Value = 15080
$ws.Range("I80").Value = 11590
$ws.Range("J80").Value = 21187.5
$ws.Range("K80").Value = 11590
$ws.Range("L80").Value = 21187.5
$ws.Range("M80").Value = -10592
$ws.Range("N80").Value = -23183.5
$ws.Range("H83").Value = 15080
$ws.Range("I83").Value = 11590
$ws.Range("J83").Value = 21187.5
$ws.Range("K83").Value = 57950
$ws.Range("L83").Value = 105937.5
$ws.Range("M83").Value = -52958
$ws.Range("N83").Value = -115921.5
$ws.Range("H97").Value = 3750
$ws.Range("J97").Value = 4500
$ws.Range("L97").Value = 4500
$ws.Range("N97").Value = -5492
$ws.Range("H126").Value = 8012
$ws.Range("I126").Value = 8012
$ws.Range("K126").Value = 24036
$ws.Range("M126").Value = -21566
$ws.Range("H132").Value = 2134.5
$ws.Range("I132").Value = 996
$ws.Range("J132").Value = 5550
$ws.Range("K132").Value = 2988
$ws.Range("L132").Value = 16650
$ws.Range("M132").Value = -458
$ws.Range("N132").Value = -21710

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 999.5
$ws.Range("I13").Value = 999.5
$ws.Range("K13").Value = 999.5
$ws.Range("M13").Value = -859.5
$ws.Range("H18").Value = 332502.25
$ws.Range("I18").Value = 332502.25
$ws.Range("K18").Value = 332502.25
$ws.Range("M18").Value = -332330.25
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H82").Value = 383.33334
$ws.Range("J82").Value = 332.5
$ws.Range("L82").Value = 332.5
$ws.Range("N82").Value = -1054.5
$ws.Range("H85").Value = 383.33334
$ws.Range("J85").Value = 332.5
$ws.Range("L85").Value = 332.5
$ws.Range("N85").Value = -2828.5
$ws.Range("H136").Value = 100989.6
$ws.Range("I136").Value = 5666.3335
$ws.Range("J136").Value = 243974.5
$ws.Range("K136").Value = 16999.0005
$ws.Range("L136").Value = 731923.5
$ws.Range("M136").Value = -14449.0005
$ws.Range("N136").Value = -737023.5
$ws.Range("H138").Value = 60001
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 60001
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 60001
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -70281

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1500
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2500
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = -328
$ws.Range("N17").Value = -2844
$ws.Range("H20").Value = 2002499.5
$ws.Range("I20").Value = 4999
$ws.Range("J20").Value = 4000000
$ws.Range("K20").Value = 4999
$ws.Range("L20").Value = 4000000
$ws.Range("N20").Value = -4000480
$ws.Range("M20").Value = -4759
$ws.Range("H54").Value = 23599.8
$ws.Range("I54").Value = 19333
$ws.Range("J54").Value = 30000
$ws.Range("K54").Value = 19333
$ws.Range("L54").Value = 30000
$ws.Range("M54").Value = -18813
$ws.Range("N54").Value = -31040
$ws.Range("H62").Value = 3586.1428
$ws.Range("I62").Value = 3220.6
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 3220.6
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -2596.6
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 3586.1428
$ws.Range("I65").Value = 3220.6
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 16103
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -12983
$ws.Range("N65").Value = -28740
$ws.Range("H107").Value = 703.7
$ws.Range("I107").Value = 703.7
$ws.Range("K107").Value = 2111.1
$ws.Range("M107").Value = -191.1000000000004
$ws.Range("H132").Value = 1243.75
$ws.Range("I132").Value = 1158.3334
$ws.Range("K132").Value = 3475.0002
$ws.Range("M132").Value = -945.0001999999999
